# Update gh-pages to output generated at 456a3b4
# Applies the "想去人数" (want-to-go count) increments, a price-to-sold-out
# change, and a data refresh (new event replacing/shifting two rows) on the
# "全部类型" aggregate sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览" (exhibitions) - "想去人数" (column F) count bumps
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value  = 1927
$ws1.Range("F8").Value  = 5297
$ws1.Range("F9").Value  = 1484
$ws1.Range("F10").Value = 152
$ws1.Range("F11").Value = 3078
$ws1.Range("F15").Value = 4218
$ws1.Range("F16").Value = 1010
$ws1.Range("F21").Value = 17
$ws1.Range("F22").Value = 129
$ws1.Range("F24").Value = 963
$ws1.Range("F29").Value = 1074
$ws1.Range("F30").Value = 361
$ws1.Range("F31").Value = 35
$ws1.Range("F32").Value = 124
$ws1.Range("F34").Value = 244
$ws1.Range("F35").Value = 1640
$ws1.Range("F36").Value = 2161
$ws1.Range("F37").Value = 1009
$ws1.Range("F40").Value = 600
$ws1.Range("F41").Value = 272
$ws1.Range("F43").Value = 648
$ws1.Range("F44").Value = 393
$ws1.Range("F45").Value = 315
$ws1.Range("F47").Value = 131

# ---------------------------------------------------------------------
# Sheet "演出" (performances) - ticket for row 17 became unavailable
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("G17").Value = "不可售"

# ---------------------------------------------------------------------
# Sheet "本地生活" (local life) - "想去人数" bump
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 729

# ---------------------------------------------------------------------
# Sheet "全部类型" (all types, aggregate) - "想去人数" count bumps
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 729
$ws4.Range("F6").Value  = 1927
$ws4.Range("F7").Value  = 5297
$ws4.Range("F8").Value  = 1484
$ws4.Range("F9").Value  = 152
$ws4.Range("F11").Value = 3078
$ws4.Range("F14").Value = 4218
$ws4.Range("F15").Value = 1010
$ws4.Range("F23").Value = 17
$ws4.Range("F26").Value = 963
$ws4.Range("F32").Value = 1074
$ws4.Range("F33").Value = 361
$ws4.Range("F34").Value = 35
$ws4.Range("F35").Value = 1640
$ws4.Range("F36").Value = 2161
$ws4.Range("F37").Value = 1009

# Row 40 is refreshed with a brand-new event (青城山下 concert); the event
# previously at row 40 ("乌托邦次元聚会3.0") slides down into row 41,
# replacing the "理查德·克莱德曼钢琴音乐会" entry that used to live there
# (that entry is dropped from this aggregate sheet; it is unaffected on the
# "演出" sheet where it already exists, just its price turned "不可售").
# (B40/B41 use a leading apostrophe so the "yyyy-MM-dd"-looking text is
# kept as plain text instead of being auto-converted to a date serial,
# matching the original file's inlineStr representation.)
$ws4.Range("B40").Value = "'2024-06-23"
$ws4.Range("C40").Value = "杭州·【早鸟5折】中西合奏·再现经典《青城山下·千年等一回》传世国风跨界音乐会"
$ws4.Range("D40").Value = "曙光路31号 浙江音乐厅"
$ws4.Range("E40").Value = "2024.06.23 15:00-06.23 21:00"
$ws4.Range("F40").Value = 1
$ws4.Range("G40").Value = 50
$ws4.Range("H40").Value = "https://show.bilibili.com/platform/detail.html?id=84597"
$ws4.Range("I40").Value = "//i2.hdslb.com/bfs/openplatform/202404/jNu5hjYv1713514034369.jpeg"

$ws4.Range("B41").Value = "'2024-06-29"
$ws4.Range("C41").Value = "杭州·乌托邦次元聚会3.0·二次元全女性夜场"
$ws4.Range("D41").Value = "保淑路2号 The Queen皇后"
$ws4.Range("E41").Value = "2024.06.29 13:00-06.29 19:00"
$ws4.Range("F41").Value = 245
$ws4.Range("G41").Value = 188
$ws4.Range("H41").Value = "https://show.bilibili.com/platform/detail.html?id=84558"
$ws4.Range("I41").Value = "//i2.hdslb.com/bfs/openplatform/202404/XyOkWYv31713435061841.jpeg"

$ws4.Range("F42").Value = 600
$ws4.Range("F43").Value = 272
$ws4.Range("F44").Value = 648
$ws4.Range("F45").Value = 393
$ws4.Range("F46").Value = 315
$ws4.Range("F48").Value = 131
